# "Generate Report for Handback" -- refresh the localization-status report
# after a handback: update statuses/dates for the handed-back files, add
# "latest target file" links, and resize the columns that now hold longer
# text.

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfecaa25e74bf8730fcc63c878f88b807976b7b8/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: the status shown for both languages moves from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------
# zh-cn sheet: record the target file + handback file + handback time,
# and point "Latest Target File" at the source markdown.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-23 10:38:16"
$wsZh.Range("K3").Value = "2016-08-23 10:38:16"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, $null, $null, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, $null, $null, "a.md") | Out-Null

$wsZh.Columns.Item(3).AutoFit()
$wsZh.Columns.Item(10).ColumnWidth = $wsZh.Columns.Item(7).ColumnWidth

# ---------------------------------------------------------------------
# de-de sheet: same refresh, but the handback timestamp differs and the
# target/handback file names carry the de-de xliff name.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-23 10:38:23"
$wsDe.Range("K3").Value = "2016-08-23 10:38:23"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, $null, $null, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, $null, $null, "a.md") | Out-Null

$wsDe.Columns.Item(3).AutoFit()
$wsDe.Columns.Item(10).ColumnWidth = $wsDe.Columns.Item(7).ColumnWidth
